$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a string value to a cell while avoiding Excel's automatic
# type-sniffing (e.g. turning "2026-02-17" into a date serial number). We set
# the literal text as a formula returning a string literal, then "flatten"
# it to a plain value via Copy + PasteSpecial(values only). This keeps the
# stored cell as a literal string with no left-over formula and no new
# number-format/style entries.
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy($cell)
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

# ===========================================================================
# Sheet "Summary"
# ===========================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.67    # Current Capital
$summary.Range("B4").Value = -0.33      # Total P&L $
$summary.Range("B5").Value = -0.25      # Total P&L %
$summary.Range("B6").Value = 26         # Total Trades
$summary.Range("B8").Value = 13         # Losing Trades
$summary.Range("B9").Value = 26.92      # Win Rate %

# ===========================================================================
# Sheet "Strategy Status" (MarketMaking row)
# ===========================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.67
$status.Range("D4").Value = 26
$status.Range("E4").Value = -0.33
$status.Range("F4").Value = -0.33
$status.Range("G4").Value = 26.92

# ===========================================================================
# New trade row (#26 -> row 27) appended to both "All Trades" and
# "MarketMaking" sheets.
# ===========================================================================
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(27, 1).Value = 26

    Set-LiteralText $ws.Cells.Item(27, 2) "2026-02-17"
    Set-LiteralText $ws.Cells.Item(27, 3) "15:19:29"
    Set-LiteralText $ws.Cells.Item(27, 4) "MarketMaking"
    Set-LiteralText $ws.Cells.Item(27, 5) "UP"

    $ws.Cells.Item(27, 6).Value = 0.84
    $ws.Cells.Item(27, 7).Value = 0.63

    Set-LiteralText $ws.Cells.Item(27, 8) "CLOSED"

    $ws.Cells.Item(27, 9).Value = -25
    $ws.Cells.Item(27, 10).Value = -0.21
    $ws.Cells.Item(27, 11).Value = 99.67
    $ws.Cells.Item(27, 12).Value = 0
    $ws.Cells.Item(27, 13).Value = 0
    $ws.Cells.Item(27, 14).Value = 0.6

    Set-LiteralText $ws.Cells.Item(27, 15) "Normal spread capture: 19600 bps"
    Set-LiteralText $ws.Cells.Item(27, 16) "early_exit"

    $ws.Cells.Item(27, 17).Value = 0.13
}
